$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 49, pushing the existing rows 49-96
# (and all their data) down to 51-98. Excel's Insert() carries the
# formatting (e.g. the date style on column D) from the row being
# displaced, which matches the target workbook's styling.
$ws.Rows.Item(49).Insert()
$ws.Rows.Item(49).Insert()

# Populate the first new row (row 49) with its data.
$ws.Range("A49").Value = 6
$ws.Range("B49").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 45126
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = 100112035
$ws.Range("G49").Value = "Bruselas (repollito)"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 150
$ws.Range("K49").Value = 16000
$ws.Range("L49").Value = 16000
$ws.Range("M49").Value = 16000
$ws.Range("N49").Value = "$/malla 15 kilos"
$ws.Range("O49").Value = "Provincia de Quillota"
$ws.Range("P49").Value = 1067
$ws.Range("Q49").Value = 15
$ws.Range("R49").Value = "Hortaliza"

# Populate the second new row (row 50) with its data.
$ws.Range("A50").Value = 6
$ws.Range("B50").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 45126
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = 100112035
$ws.Range("G50").Value = "Bruselas (repollito)"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 580
$ws.Range("K50").Value = 17000
$ws.Range("L50").Value = 18000
$ws.Range("M50").Value = 17397
$ws.Range("N50").Value = "$/malla 15 kilos"
$ws.Range("O50").Value = "Provincia de Quillota"
$ws.Range("P50").Value = 1160
$ws.Range("Q50").Value = 15
$ws.Range("R50").Value = "Hortaliza"

# Ensure the date columns use the same number format as the rest of
# column D (matches style index 2 used throughout the sheet).
$ws.Range("D49").NumberFormat = $ws.Range("D51").NumberFormat()
$ws.Range("D50").NumberFormat = $ws.Range("D51").NumberFormat()
